$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Academic resource entry from CURCE to ICA
$ws.Range("A2").Value = "Aca_ICA"
$ws.Range("B2").Value = "International Crossroads Albany (ICA)"
$ws.Range("D2").Value = "International Crossroads Albany (ICA)"

# Update the static generated-code text cells (copies of the B4:B6 formula
# results) to reflect the new ICA entry instead of the old CURCE one
$ws.Range("A8").Value = $ws.Range("B4").Value()
$ws.Range("A9").Value = $ws.Range("B5").Value()
$ws.Range("A10").Value = $ws.Range("B6").Value()

# Update the selection to match the new active cell / range on the sheet
$ws.Range("A8:A10").Select()
